{"js": "// Surat Keterangan Tidak Mampu \u2014 field updates (nama, NIK, jenis kelamin,\n// tempat/tgl lahir, beasiswa keterangan, tanggal surat).\n//\n// Each target string is unique in the document body, so a scoped,\n// case-sensitive search-and-replace is unambiguous for every field.\n\nconst replacements = [\n  [\"JULAEHA\", \"KIRMAN\"],\n  [\"3208074101710002\", \"3208270107850249\"],\n  [\"Perempuan\", \"Laki-laki\"],\n  [\"KUNINGAN, 25934\", \"KUNINGAN, 31229\"],\n  [\"beasiswa\", \"z\"],\n  [\"Ciawigebang, 02 Oktober 2017\", \"Ciawigebang, 16 Oktober 2017\"],\n];\n\nconst body = context.document.body;\n\nfor (const [searchText, newText] of replacements) {\n  const results = body.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Surat Keterangan Tidak Mampu - field updates (nama, NIK, jenis kelamin,\n# tempat/tgl lahir, beasiswa keterangan, tanggal surat).\n#\n# Each target string is unique in the document body, so a scoped,\n# case-sensitive Find/Replace (ReplaceOne) is unambiguous for every field.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n}\n\nReplace-Text \"JULAEHA\" \"KIRMAN\"\nReplace-Text \"3208074101710002\" \"3208270107850249\"\nReplace-Text \"Perempuan\" \"Laki-laki\"\nReplace-Text \"KUNINGAN, 25934\" \"KUNINGAN, 31229\"\nReplace-Text \"beasiswa\" \"z\"\nReplace-Text \"Ciawigebang, 02 Oktober 2017\" \"Ciawigebang, 16 Oktober 2017\"\n"}
